$wb = $excel.ActiveWorkbook

# --- GraphChi Implementation Notes: insert a new note row (row 33) ---
$ws4 = $wb.Worksheets.Item("GraphChi Implementation Notes")

# Insert a new blank row at position 33, pushing old rows 33-39 down to 34-40.
[void]$ws4.Rows.Item(33).Insert()

# The row-insert shift auto-populated cells in A:D from the row above;
# clear them so only column E holds content, matching the target layout.
$ws4.Range("A33:D33").Clear()

# Copy the formatting from the similarly-styled note cell above (E31/E32,
# style index 34) onto the new note cell, then set its text and row height.
$ws4.Range("E31").Copy()
$ws4.Range("E33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("E33").Value = "8/17/2015 2:38 PM`n!!! Need to look into friendsoffriends.java. Seems useful for DTC"
$ws4.Rows.Item(33).RowHeight = 43.2

# --- Resources and GraphChi Info: move selection ---
$ws3 = $wb.Worksheets.Item("Resources and GraphChi Info")
[void]$ws3.Range("D12").Select()

# --- GraphChi Implementation Notes becomes the active/selected sheet ---
[void]$ws4.Activate()
[void]$ws4.Range("E34").Select()

Write-Output "done"
